$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Numero_de_obligacion"
$ws.Range("I1").Value = "Fecha_de_Obligacion"
$ws.Range("M1").Value = "Identificacion_Beneficiario"
